# LOB1052.xlsx update
# - "Objetivos:" (row 10) gets its real Portuguese objectives text (it previously,
#   erroneously, held the professor's name).
# - A new row is inserted for "Docentes responsáveis:" (row 12) so the professor's
#   name has its own B/C value row (row 13).
# - "Programa resumido:" / "Programa:" / "Método:" / "Critério:" / "Norma de
#   recuperação:" / "Bibliografia:" all get their correct content; previously each
#   held a value that actually belonged to a different field (an off-by-one content
#   shift), and "Bibliografia:" had no content at all.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix "Objetivos:" value row (row 10) -----------------------------------
$ws.Range("B10").Value = "Apresentar os conceitos teóricos e aplicações da Integração de funções de varias variáveis reais e o calculo vetorial."
$ws.Range("C10").Value = "Apresentar os conceitos teóricos e aplicações da Integração de funções de varias variáveis reais e o calculo vetorial."

# --- Insert the missing "Docentes responsáveis:" value row ----------------
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Clear()

# Copy the B/C number/alignment/font formatting down from the row below
# (style 2 for B, style 3 for C) before filling in the values.
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B13").Value = "6270264 - Juan Fernando Zapata Zapata"
$ws.Range("C13").Value = "6270264 - Juan Fernando Zapata Zapata"

# --- Fix "Programa resumido:" value row (now row 14) -----------------------
$ws.Range("B14").Value = "Integrais Múltiplas, Integrais de Linha, Integrais de Superfície e Teorema de Stokes."
$ws.Range("C14").Value = "Integrais Múltiplas, Integrais de Linha, Integrais de Superfície e Teorema de Stokes."

# --- Fix "Programa:" value row (now row 16) ---------------------------------
$programa = "•Integrais Múltiplas:  Integrais Duplas e triplas, integrais iteradas e o Teorema de Fubinni, teorema de mudança de variáveis, Aplicações.•Campos de vetores: Definição, Operadores rotacional e divergente para campos de vetores. •Integral de Linha: Definição, trabalho e energia, Teorema fundamental da integral de linha, Campos conservativos, teorema de Green, Fluxo de um campo de vetores sobre uma curva.•Integrais de superfície: Superfícies parametrizadas, orientação de superfícies, Integrais de Superfície e aplicações.•Teoremas Vetoriais: Teorema de Stokes e Teorema da divergência, lei de indução de Faraday e equação de continuidade dos fluidos."
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

# --- Fix "Método:" value row (now row 19) -----------------------------------
$metodo = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# --- Fix "Critério:" value row (now row 20) ---------------------------------
$ws.Range("B20").Value = "NF≥ 5,0."
$ws.Range("C20").Value = "NF≥ 5,0."

# --- Fix "Norma de recuperação:" value row (now row 21) ---------------------
$norma = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
$ws.Range("B21").Value = $norma
$ws.Range("C21").Value = $norma

# --- Fix "Bibliografia:" value row (now row 22) -----------------------------
$biblio = "1. H. L. Guidorizzi, UM CURSO DE CÁLCULO, volume III. Livros Técnicos e Científicos, Rio de Janeiro.2. W. Kaplan, CÁLCULO AVANÇADO, volume I, Edgard Blücher, 1972.3. Stewart, CÁLCULO, volume II, Editora Pioneira Thomson Leaming.4.BUSS, Mirian ; FLEMMING, Diva Marília. Calculo B. 2. ed. São Paulo:Pearson, 2007."
$ws.Range("B22").Value = $biblio
$ws.Range("C22").Value = $biblio
